$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5002512.5
$ws.Range("J17").Value = 5002512.5
$ws.Range("L17").Value = 15007537.5
$ws.Range("N17").Value = -15007873.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4380.2085
$ws.Range("I32").Value = 4380.2085
$ws.Range("K32").Value = 4380.2085
$ws.Range("M32").Value = -4093.2085

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 27400
$ws.Range("I61").Value = 31750
$ws.Range("K61").Value = 31750
$ws.Range("M61").Value = -31538

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5497
$ws.Range("I74").Value = 3954.5833
$ws.Range("K74").Value = 3954.5833
$ws.Range("M74").Value = -3080.5833

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5497
$ws.Range("I77").Value = 3954.5833
$ws.Range("K77").Value = 19772.9165
$ws.Range("M77").Value = -15404.9165

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 27400
$ws.Range("I136").Value = 31750
$ws.Range("K136").Value = 95250
$ws.Range("M136").Value = -92700

# BSM row 25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 675.6667

# BSM row 55
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 74499.664
$ws.Range("J55").Value = 74499.664
$ws.Range("L55").Value = 74499.664
$ws.Range("N55").Value = -75045.664

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2263.5806
$ws.Range("I105").Value = 2049.16
$ws.Range("K105").Value = 2049.16
$ws.Range("M105").Value = -302.1599999999999

# BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 64999
$ws.Range("J126").Value = 64999
$ws.Range("L126").Value = 64999
$ws.Range("N126").Value = -74879

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4307.3335
$ws.Range("I134").Value = 3812.5454
$ws.Range("K134").Value = 11437.6362
$ws.Range("M134").Value = -8902.636200000001

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 78995
$ws.Range("J140").Value = 78995
$ws.Range("L140").Value = 78995
$ws.Range("N140").Value = -89355

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 77.22221999999999
$ws.Range("I7").Value = 81.17646999999999
$ws.Range("K7").Value = 81.17646999999999
$ws.Range("M7").Value = 31.82353000000001

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30308158
$ws.Range("J31").Value = 7786.533
$ws.Range("L31").Value = 7786.533
$ws.Range("N31").Value = -8376.532999999999

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 30308158
$ws.Range("J34").Value = 7786.533
$ws.Range("L34").Value = 7786.533
$ws.Range("N34").Value = -8190.533

# CRP row 48
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8379.393
$ws.Range("I58").Value = 5116.5
$ws.Range("K58").Value = 5116.5
$ws.Range("M58").Value = -4913.5

# CRP row 97
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 25263.5
$ws.Range("J97").Value = 26019.715
$ws.Range("L97").Value = 26019.715
$ws.Range("N97").Value = -28001.715

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 641.8421
$ws.Range("I105").Value = 647.5
$ws.Range("K105").Value = 647.5
$ws.Range("M105").Value = 1099.5

# CRP row 110
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 72124
$ws.Range("J110").Value = 79498.664
$ws.Range("L110").Value = 79498.664
$ws.Range("N110").Value = -87678.664

# CRP row 115
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H115").Value = 49999.668
$ws.Range("I115").Value = 49999
$ws.Range("K115").Value = 49999
$ws.Range("M115").Value = -48824

# CRP row 129
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 80000
$ws.Range("J129").Value = 80000
$ws.Range("L129").Value = 80000
$ws.Range("N129").Value = -90000

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4819.4287
$ws.Range("I134").Value = 2398.3333
$ws.Range("J134").Value = 9177.4
$ws.Range("K134").Value = 7194.999899999999
$ws.Range("L134").Value = 27532.2
$ws.Range("M134").Value = -4659.999899999999
$ws.Range("N134").Value = -32602.2

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8379.393
$ws.Range("I136").Value = 5116.5
$ws.Range("K136").Value = 15349.5
$ws.Range("M136").Value = -12799.5

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 63673588
$ws.Range("I4").Value = 83217930
$ws.Range("K4").Value = 249653790
$ws.Range("M4").Value = -249653678

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 558.8889
$ws.Range("I5").Value = 447.85715
$ws.Range("K5").Value = 1343.57145
$ws.Range("M5").Value = -1231.57145

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 527.05884
$ws.Range("I12").Value = 34.714287
$ws.Range("K12").Value = 104.142861
$ws.Range("M12").Value = 68.857139

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 500
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 500
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -1416
$ws.Range("N34").Value = -1668

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1437.6

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 425.33334
$ws.Range("I46").Value = 250.66667
$ws.Range("J46").Value = 600
$ws.Range("K46").Value = 752.00001
$ws.Range("L46").Value = 1800
$ws.Range("M46").Value = -661.00001
$ws.Range("N46").Value = -1982

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4219.4443
$ws.Range("I55").Value = 995
$ws.Range("J55").Value = 6799
$ws.Range("K55").Value = 2985
$ws.Range("L55").Value = 20397
$ws.Range("M55").Value = -2808
$ws.Range("N55").Value = -20751

# CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 41749.875
$ws.Range("I56").Value = 41749.875
$ws.Range("K56").Value = 41749.875
$ws.Range("M56").Value = -41219.875

# CUL row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 13149
$ws.Range("J69").Value = 13149
$ws.Range("L69").Value = 39447
$ws.Range("N69").Value = -41069

# CUL row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 13149
$ws.Range("J72").Value = 13149
$ws.Range("L72").Value = 118341
$ws.Range("N72").Value = -126453

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 303318.3
$ws.Range("J107").Value = 588546.3
$ws.Range("L107").Value = 1765638.9
$ws.Range("N107").Value = -1769478.9

# CUL row 119
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 366.66666
$ws.Range("I119").Value = 366.66666
$ws.Range("K119").Value = 1099.99998
$ws.Range("M119").Value = 3738.00002

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 558.8889
$ws.Range("I135").Value = 447.85715
$ws.Range("K135").Value = 4030.71435
$ws.Range("M135").Value = -1495.71435

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5401.1816
$ws.Range("I137").Value = 3492.625
$ws.Range("J137").Value = 10490.667
$ws.Range("K137").Value = 10477.875
$ws.Range("L137").Value = 31472.001
$ws.Range("M137").Value = -5377.875
$ws.Range("N137").Value = -41672.001

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7236.6875
$ws.Range("I46").Value = 7236.6875
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 7236.6875
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -7048.6875
$ws.Range("N46").ClearContents()

# LTW row 74
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 54293.332
$ws.Range("I74").Value = 54293.332
$ws.Range("K74").Value = 54293.332
$ws.Range("M74").Value = -53295.332

# LTW row 77
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 54293.332
$ws.Range("I77").Value = 54293.332
$ws.Range("K77").Value = 162879.996
$ws.Range("M77").Value = -157887.996

# LTW row 108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 74996
$ws.Range("J108").Value = 74996
$ws.Range("L108").Value = 74996
$ws.Range("N108").Value = -82676

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 111119920
$ws.Range("I136").Value = 76932216
$ws.Range("K136").Value = 230796648
$ws.Range("M136").Value = -230794098

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 50017.57
$ws.Range("J46").Value = 50017.57
$ws.Range("L46").Value = 50017.57
$ws.Range("N46").Value = -50479.57

# WVR row 52
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 30466.666
$ws.Range("J52").Value = 31400
$ws.Range("L52").Value = 31400
$ws.Range("N52").Value = -31852

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4527.175
$ws.Range("I132").Value = 3737.9355
$ws.Range("J132").Value = 7245.6665
$ws.Range("K132").Value = 11213.8065
$ws.Range("L132").Value = 21736.9995
$ws.Range("M132").Value = -8683.806500000001
$ws.Range("N132").Value = -26796.9995

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 50017.57
$ws.Range("J134").Value = 50017.57
$ws.Range("L134").Value = 150052.71
$ws.Range("N134").Value = -155122.71
